# Widen main content placeholder
#
# The "Content Placeholder" shape lives on the "Title and Content" slide
# layout (ppt/slideLayouts/slideLayout1.xml) that belongs to the deck's
# (only) slide master. There are no slides in this template presentation,
# so the edit is made directly on the layout placeholder's geometry.

$p = $ppt.ActivePresentation

# --- Widen the main content placeholder on the "Title and Content" layout ---
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

$contentLayout = $null
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $candidate = $master.CustomLayouts.Item($i)
    if ($candidate.Name -eq "Title and Content") {
        $contentLayout = $candidate
    }
}

$contentPh = $null
for ($i = 1; $i -le $contentLayout.Shapes.Count; $i++) {
    $shp = $contentLayout.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder") {
        $contentPh = $shp
    }
}

# New geometry, expressed in points (1 pt = 12700 EMU) so that the
# round-tripped Single-precision COM values serialize back to the exact
# target EMU: off x=167054 y=923544, ext cx=11887200 cy=5733288.
$contentPh.Left = 13.1539
$contentPh.Top = 72.72
$contentPh.Width = 936.0
$contentPh.Height = 451.44

# --- Refresh the cached "updates automatically" date fields ---
# (slide master + "Title Slide" layout Date placeholders) so they read the
# same new date, mirroring what happens when PowerPoint re-caches an
# auto date field on save.
$newDateText = "4/24/2022"

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -eq "Date Placeholder 3") {
        $shp.TextFrame.TextRange.Text = $newDateText
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -eq "Date Placeholder 3" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -ne "") {
                $shp.TextFrame.TextRange.Text = $newDateText
            }
        }
    }
}
